# Apply weekly refresh of Fruta/Hortaliza (Granada) data.
# The diff shows that the values in columns D and L:T for rows 2-15 are
# simply reshuffled among the existing rows (same set of records, new
# row order / new "current" assignment). Capture the original values for
# each row first, then write them back out to their new destination row
# per the mapping below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (source row values, as they
# existed in the original workbook, should end up in destination row).
$map = @{
    2  = 7
    3  = 9
    4  = 8
    5  = 2
    6  = 14
    7  = 3
    8  = 4
    9  = 10
    10 = 5
    11 = 6
    12 = 13
    13 = 15
    14 = 12
    15 = 11
}

# Columns whose values move together with the row (D and L..T).
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot original values for every source row/column before any writes.
# Value2 is used (rather than Value) so that dates come back as raw
# numeric serials instead of formatted date strings.
$orig = @{}
foreach ($r in 2..15) {
    $orig[$r] = @{}
    foreach ($col in $cols) {
        $orig[$r][$col] = $ws.Range("$col$r").Value2()
    }
}

# Write the snapshot values into their new destination rows.
foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $orig[$srcRow][$col]
    }
}
